$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the text values that changed (shared strings reused elsewhere are
# updated in place for these specific cells).
$ws.Range("H5").Value = "Only After 7:00pm"
$ws.Range("H6").Value = "no 1/10/14-1/13/14, no 1/6/14"
$ws.Range("H7").Value = "No 1/12/14"

# Update the active cell selection to H5 on this sheet.
$ws.Activate()
$ws.Range("H5").Select()
